$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.7964671443145639
$ws.Range("D2").Value = 0.7906774682135947
$ws.Range("E2").Value = 0.7814381039512749
$ws.Range("F2").Value = 0.777841828637237
$ws.Range("G2").Value = 0.7922867964911615
$ws.Range("H2").Value = 0.8402330979296048
$ws.Range("I2").Value = 0.94799791016554
$ws.Range("J2").Value = 1.3008905245729487
$ws.Range("K2").Value = "#NUM!"
$ws.Range("L2").Value = "#NUM!"
$ws.Range("C3").Value = 0.6628295253534962
$ws.Range("D3").Value = 0.663006218440504
$ws.Range("E3").Value = 0.6630935352519158
$ws.Range("F3").Value = 0.6639640597683897
$ws.Range("G3").Value = 0.6679471918592815
$ws.Range("H3").Value = 0.679575151038384
$ws.Range("I3").Value = 0.7058480046218621
$ws.Range("J3").Value = 0.7097351459341179
$ws.Range("K3").Value = "#NUM!"
$ws.Range("L3").Value = "#NUM!"
$ws.Range("C4").Value = 0.5099930956620913
$ws.Range("D4").Value = 0.45711599013229465
$ws.Range("E4").Value = 0.40607704333716504
$ws.Range("F4").Value = 0.36470317651973394
$ws.Range("G4").Value = 0.3406899078713379
$ws.Range("H4").Value = 0.3404923255650547
$ws.Range("I4").Value = 0.3716256617972899
$ws.Range("J4").Value = 0.49627489864770086
$ws.Range("K4").Value = "#NUM!"
$ws.Range("L4").Value = "#NUM!"
$ws.Range("C5").Value = 0.3020070628335136
$ws.Range("D5").Value = 0.3542523085794068
$ws.Range("E5").Value = 0.4075166639982665
$ws.Range("F5").Value = 0.4615750624469502
$ws.Range("G5").Value = 0.515170925893117
$ws.Range("H5").Value = 0.5643871840916832
$ws.Range("I5").Value = 0.5984377770153521
$ws.Range("J5").Value = 0.6588468883788751
$ws.Range("K5").Value = "#NUM!"
$ws.Range("L5").Value = "#NUM!"
$ws.Range("C6").Value = 0.8392883244104475
$ws.Range("D6").Value = 0.8280065221462679
$ws.Range("E6").Value = 0.8170804974956407
$ws.Range("F6").Value = 0.8065108131163876
$ws.Range("G6").Value = 0.797218211420159
$ws.Range("H6").Value = 0.7914882516171455
$ws.Range("I6").Value = 0.7926615351408527
$ws.Range("J6").Value = 0.8105991436635933
$ws.Range("K6").Value = "#NUM!"
$ws.Range("L6").Value = "#NUM!"
$ws.Range("C7").Value = 0.8264785300783722
$ws.Range("D7").Value = 0.8245404310007802
$ws.Range("E7").Value = 0.8283598582102396
$ws.Range("F7").Value = 0.8386558104572086
$ws.Range("G7").Value = 0.8561624135773356
$ws.Range("H7").Value = 0.8817507540953591
$ws.Range("I7").Value = 0.9162287135217525
$ws.Range("J7").Value = 0.9589327580386309
$ws.Range("K7").Value = "#NUM!"
$ws.Range("L7").Value = "#NUM!"
$ws.Range("C8").Value = -0.09997666402937262
$ws.Range("D8").Value = -0.1999575110026547
$ws.Range("E8").Value = -0.2999219231562071
$ws.Range("F8").Value = -0.3998419320349996
$ws.Range("G8").Value = -0.49969315389369456
$ws.Range("H8").Value = -0.5994563648821072
$ws.Range("I8").Value = -0.6990186661349361
$ws.Range("J8").Value = -0.7977693459062928
$ws.Range("K8").Value = "#NUM!"
$ws.Range("L8").Value = "#NUM!"
$ws.Range("C9").Value = 0.7840219905118377
$ws.Range("D9").Value = 0.8656586752333428
$ws.Range("E9").Value = 0.942982931196805
$ws.Range("F9").Value = 1.0141398597649185
$ws.Range("G9").Value = 1.0778066029841034
$ws.Range("H9").Value = 1.132994859131851
$ws.Range("I9").Value = 1.18808449546574
$ws.Range("J9").Value = 1.3635740929948894
$ws.Range("K9").Value = "#NUM!"
$ws.Range("L9").Value = "#NUM!"
$ws.Range("C10").Value = 0.5855114907477059
$ws.Range("D10").Value = 0.6707727835243569
$ws.Range("E10").Value = 0.7560569160663637
$ws.Range("F10").Value = 0.8418217979794852
$ws.Range("G10").Value = 0.9287355349196833
$ws.Range("H10").Value = 1.017706644673415
$ws.Range("I10").Value = 1.109573490104472
$ws.Range("J10").Value = 1.200207450498878
$ws.Range("K10").Value = "#NUM!"
$ws.Range("L10").Value = "#NUM!"
$ws.Range("C11").Value = 0.49595820946638425
$ws.Range("D11").Value = 0.49674148198757523
$ws.Range("E11").Value = 0.5037108304650835
$ws.Range("F11").Value = 0.5178455560495503
$ws.Range("G11").Value = 0.539624519466707
$ws.Range("H11").Value = 0.5688239232760365
$ws.Range("I11").Value = 0.6041300036879272
$ws.Range("J11").Value = 0.6396590251383449
$ws.Range("K11").Value = "#NUM!"
$ws.Range("L11").Value = "#NUM!"
